$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2-8
# from serial date 45184 (2023-09-15) to 45185 (2023-09-16)
$ws.Range("C2:C8").Value = 45185
